$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for price cells that would otherwise be
# auto-converted to numbers by Excel (values like "595.44").
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D18", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D31", "D33", "D34", "D35", "D37", "D39", "D40", "D43", "D44", "D46", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.271.57"
$ws.Range("E2").Value = "  -4.83%  "
$ws.Range("D3").Value = "3.257.70"
$ws.Range("E3").Value = "  -7.45%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "595.44"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").Value = "151.14"
$ws.Range("E6").Value = "  -12.66%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.249.30"
$ws.Range("E8").Value = "  -7.54%  "
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  -11.16%  "
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  -13.26%  "
$ws.Range("D11").Value = "6.76"
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  -12.90%  "
$ws.Range("D13").Value = "38.05"
$ws.Range("E13").Value = "  -17.98%  "
$ws.Range("D14").Value = "0.0000243"
$ws.Range("E14").Value = "  -12.06%  "
$ws.Range("D15").Value = "3.780.85"
$ws.Range("E15").Value = "  -7.56%  "
$ws.Range("D16").Value = "67.262.88"
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").Value = "3.255.81"
$ws.Range("E17").Value = "  -7.53%  "
$ws.Range("D18").Value = "545.48"
$ws.Range("E18").Value = "  -10.37%  "
$ws.Range("E19").Value = "  -6.05%  "
$ws.Range("E20").Value = "  -13.58%  "
$ws.Range("D21").Value = "15.12"
$ws.Range("E21").Value = "  -14.65%  "
$ws.Range("D22").Value = "0.763"
$ws.Range("E22").Value = "  -13.47%  "
$ws.Range("D23").Value = "7.86"
$ws.Range("E23").Value = "  -14.24%  "
$ws.Range("D24").Value = "85.51"
$ws.Range("E24").Value = "  -12.80%  "
$ws.Range("D25").Value = "13.56"
$ws.Range("E25").Value = "  -13.35%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "3.24"
$ws.Range("E27").Value = "  -13.16%  "
$ws.Range("D28").Value = "29.38"
$ws.Range("E28").Value = "  -12.82%  "
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  -11.14%  "
$ws.Range("E30").Value = "  -17.44%  "
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  -11.04%  "
$ws.Range("E32").Value = "  -12.03%  "
$ws.Range("D33").Value = "6.66"
$ws.Range("E33").Value = "  -17.66%  "
$ws.Range("D34").Value = "537.65"
$ws.Range("E34").Value = "  -15.83%  "
$ws.Range("D35").Value = "5.72"
$ws.Range("E35").Value = "  -15.65%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "0.0449"
$ws.Range("E37").Value = "  -7.36%  "
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("D39").Value = "0.0853"
$ws.Range("E39").Value = "  -14.58%  "
$ws.Range("D40").Value = "9.15"
$ws.Range("E40").Value = "  -15.25%  "
$ws.Range("E41").Value = "  -10.79%  "
$ws.Range("D42").Value = "2.930.13"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  -21.87%  "
$ws.Range("D44").Value = "0.262"
$ws.Range("E44").Value = "  -16.15%  "
$ws.Range("D45").Value = "0.0₃0582"
$ws.Range("E45").Value = "  -18.95%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  -14.86%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "26.42"
$ws.Range("E47").Value = "  -16.81%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "127.94"
$ws.Range("E49").Value = "  -4.91%  "
$ws.Range("E50").Value = "  -19.67%  "
$ws.Range("D51").Value = "0.114"
$ws.Range("E51").Value = "  -12.69%  "
